# Scheduled refresh of market-board snapshot values (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the per-class Leve-profit sheets.
# Values below come from the latest pull; cells with no meaningful
# profit figure (denominator 0 / not for sale) are cleared instead of
# left at a stale number.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2079.6428
$ws.Range("I40").Value = 1819.8
$ws.Range("J40").Value = 2224
$ws.Range("K40").Value = 1819.8
$ws.Range("L40").Value = 2224
$ws.Range("M40").Value = -1644.8
$ws.Range("N40").Value = -2574

$ws.Range("H100").Value = 2066.5
$ws.Range("I100").Value = 2066.5
$ws.Range("K100").Value = 2066.5
$ws.Range("M100").Value = -1525.5

$ws.Range("H138").Value = 10003671
$ws.Range("J138").Value = 3566.4285
$ws.Range("L138").Value = 10699.2855
$ws.Range("N138").Value = -20979.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 19694.2
$ws.Range("I28").Value = 19694.2
$ws.Range("K28").Value = 19694.2
$ws.Range("M28").Value = -19502.2

$ws.Range("H45").Value = 3900.75
$ws.Range("I45").Value = 2636.875
$ws.Range("K45").Value = 2636.875
$ws.Range("M45").Value = -2259.875

$ws.Range("H51").Value = 20000
$ws.Range("I51").Value = 20000
$ws.Range("K51").Value = 20000
$ws.Range("M51").Value = -19244

$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

$ws.Range("H53").Value = 40000
$ws.Range("I53").Value = 40000
$ws.Range("K53").Value = 40000
$ws.Range("M53").Value = -39318

$ws.Range("H99").Value = 19694.2
$ws.Range("I99").Value = 19694.2
$ws.Range("K99").Value = 19694.2
$ws.Range("M99").Value = -16699.2

$ws.Range("H122").Value = 2428.75
$ws.Range("I122").Value = 2400.4
$ws.Range("J122").Value = 2476
$ws.Range("K122").Value = 7201.200000000001
$ws.Range("L122").Value = 7428
$ws.Range("M122").Value = -4751.200000000001
$ws.Range("N122").Value = -12328

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 10152
$ws.Range("I26").Value = 10152
$ws.Range("K26").Value = 10152
$ws.Range("M26").Value = -9860

$ws.Range("H105").Value = 1874.125
$ws.Range("I105").Value = 1784.7142
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 1784.7142
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -37.71419999999989
$ws.Range("N105").Value = -5994

$ws.Range("H134").Value = 5500.4287
$ws.Range("I134").Value = 2375.75
$ws.Range("J134").Value = 9666.666999999999
$ws.Range("K134").Value = 7127.25
$ws.Range("L134").Value = 29000.001
$ws.Range("M134").Value = -4592.25
$ws.Range("N134").Value = -34070.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 252.11111
$ws.Range("I7").Value = 313.14285
$ws.Range("K7").Value = 313.14285
$ws.Range("M7").Value = -200.14285

$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H51").Value = 5995
$ws.Range("I51").Value = 5995
$ws.Range("K51").Value = 5995
$ws.Range("M51").Value = -5259

$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").ClearContents()

$ws.Range("H59").Value = 31567
$ws.Range("J59").Value = 41245
$ws.Range("L59").Value = 41245
$ws.Range("N59").Value = -43535

$ws.Range("H61").Value = 5995
$ws.Range("I61").Value = 5995
$ws.Range("K61").Value = 5995
$ws.Range("M61").Value = -5647

$ws.Range("H105").Value = 1936.1666
$ws.Range("I105").Value = 1936.1666
$ws.Range("K105").Value = 1936.1666
$ws.Range("M105").Value = -189.1666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2399.6667
$ws.Range("I109").Value = 2399.6667
$ws.Range("K109").Value = 7199.000100000001
$ws.Range("M109").Value = -6159.000100000001

$ws.Range("H139").Value = 15000
$ws.Range("I139").Value = 15000
$ws.Range("K139").Value = 45000
$ws.Range("M139").Value = -39860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 39.9
$ws.Range("I2").Value = 15.916667
$ws.Range("J2").Value = 75.875
$ws.Range("K2").Value = 15.916667
$ws.Range("L2").Value = 75.875
$ws.Range("M2").Value = 97.083333
$ws.Range("N2").Value = -301.875

$ws.Range("H70").Value = 3479.8
$ws.Range("I70").Value = 3474.75
$ws.Range("K70").Value = 3474.75
$ws.Range("M70").Value = -3204.75

$ws.Range("H73").Value = 3479.8
$ws.Range("I73").Value = 3474.75
$ws.Range("K73").Value = 3474.75
$ws.Range("M73").Value = -2538.75

$ws.Range("H99").Value = 7200
$ws.Range("I99").Value = 7200
$ws.Range("K99").Value = 7200
$ws.Range("M99").Value = -4954

$ws.Range("H102").Value = 4071.111
$ws.Range("I102").Value = 2980.8572
$ws.Range("J102").Value = 7887
$ws.Range("K102").Value = 2980.8572
$ws.Range("L102").Value = 7887
$ws.Range("M102").Value = -1358.8572
$ws.Range("N102").Value = -11131

$ws.Range("H107").Value = 395
$ws.Range("I107").Value = 293.75
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 293.75
$ws.Range("L107").Value = 800
$ws.Range("M107").Value = 1626.25
$ws.Range("N107").Value = -4640

$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H126").Value = 2287.5
$ws.Range("J126").Value = 150
$ws.Range("L126").Value = 450
$ws.Range("N126").Value = -5390

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3255.6667
$ws.Range("I46").Value = 3466.8333
$ws.Range("J46").Value = 2833.3333
$ws.Range("K46").Value = 3466.8333
$ws.Range("L46").Value = 2833.3333
$ws.Range("M46").Value = -3278.8333
$ws.Range("N46").Value = -3209.3333

$ws.Range("H122").Value = 3322.7
$ws.Range("J122").Value = 3286.3333
$ws.Range("L122").Value = 9858.999899999999
$ws.Range("N122").Value = -14758.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H136").Value = 11293.8
$ws.Range("J136").Value = 15750
$ws.Range("L136").Value = 47250
$ws.Range("N136").Value = -52350
